# Insert a new data row at row 108 (pushing the existing rows 108..181 down
# to 109..182, matching the target diff) and populate it with the new
# weekly record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("108:108").Insert()

$ws.Range("A108").Value = 11
$ws.Range("B108").Value = "Vega Monumental Concepción"
$ws.Range("C108").Value = "Bíobío"
$ws.Range("D108").Value = 44957
$ws.Range("E108").Value = 8
$ws.Range("F108").Value = 100112043
$ws.Range("G108").Value = "Pepino ensalada"
$ws.Range("H108").Value = "Sin especificar"
$ws.Range("I108").Value = "Primera"
$ws.Range("J108").Value = 110
$ws.Range("K108").Value = 6500
$ws.Range("L108").Value = 7000
$ws.Range("M108").Value = 6727
$ws.Range("N108").Value = "`$/caja 60 unidades"
$ws.Range("O108").Value = "Región de Arica y Parinacota"
$ws.Range("P108").Value = 112
$ws.Range("Q108").Value = 60
$ws.Range("R108").Value = "Hortaliza"
